$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 598.5
$ws.Range("I2").Value = 464.66666
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 464.66666
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -351.66666
$ws.Range("N2").Value = -1226
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H33").Value = 411.70587
$ws.Range("I33").Value = 411.70587
$ws.Range("K33").Value = 411.70587
$ws.Range("M33").Value = -182.70587
$ws.Range("I40").Value = 1699.1666
$ws.Range("J40").Value = 2799.6667
$ws.Range("K40").Value = 1699.1666
$ws.Range("L40").Value = 2799.6667
$ws.Range("M40").Value = -1524.1666
$ws.Range("N40").Value = -3149.6667
$ws.Range("H53").Value = 500
$ws.Range("I53").Value = 500
$ws.Range("K53").Value = 500
$ws.Range("M53").Value = 137
$ws.Range("H70").Value = 2499.75
$ws.Range("J70").Value = 2499.75
$ws.Range("L70").Value = 7499.25
$ws.Range("N70").Value = -8039.25
$ws.Range("H73").Value = 2499.75
$ws.Range("J73").Value = 2499.75
$ws.Range("L73").Value = 7499.25
$ws.Range("N73").Value = -9371.25
$ws.Range("H80").Value = 995
$ws.Range("J80").Value = 995
$ws.Range("L80").Value = 2985
$ws.Range("N80").Value = -4981
$ws.Range("H83").Value = 995
$ws.Range("J83").Value = 995
$ws.Range("L83").Value = 8955
$ws.Range("N83").Value = -18939
$ws.Range("H107").Value = 56805.312
$ws.Range("I107").Value = 75607
$ws.Range("K107").Value = 75607
$ws.Range("M107").Value = -73687
$ws.Range("H127").Value = 951.4
$ws.Range("I127").Value = 837.5714
$ws.Range("K127").Value = 2512.7142
$ws.Range("M127").Value = 2447.2858
$ws.Range("H137").Value = 583
$ws.Range("J137").Value = 666
$ws.Range("L137").Value = 1998
$ws.Range("N137").Value = -7098

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2199.8333
$ws.Range("J45").Value = 1925
$ws.Range("L45").Value = 1925
$ws.Range("N45").Value = -2679
$ws.Range("H61").Value = 1918
$ws.Range("I61").Value = 1918
$ws.Range("K61").Value = 1918
$ws.Range("M61").Value = -1706
$ws.Range("H132").Value = 2669.5715
$ws.Range("I132").Value = 2447.8333
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7343.499899999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4813.499899999999
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 1918
$ws.Range("I136").Value = 1918
$ws.Range("K136").Value = 5754
$ws.Range("M136").Value = -3204

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 28876.5
$ws.Range("J95").Value = 28876.5
$ws.Range("L95").Value = 28876.5
$ws.Range("N95").Value = -34368.5
$ws.Range("H105").Value = 1990.3334
$ws.Range("I105").Value = 1990.3334
$ws.Range("K105").Value = 1990.3334
$ws.Range("M105").Value = -243.3334
$ws.Range("H107").Value = 53223.625
$ws.Range("I107").Value = 69683.336
$ws.Range("K107").Value = 69683.336
$ws.Range("M107").Value = -67763.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1010.7778
$ws.Range("I22").Value = 916.1667
$ws.Range("K22").Value = 916.1667
$ws.Range("M22").Value = -566.1667
$ws.Range("H107").Value = 670.1667
$ws.Range("I107").Value = 473.66666
$ws.Range("K107").Value = 473.66666
$ws.Range("M107").Value = 1446.33334
$ws.Range("H124").Value = 48469
$ws.Range("J124").Value = 48469
$ws.Range("L124").Value = 48469
$ws.Range("N124").Value = -53379
$ws.Range("H132").Value = 2023.125
$ws.Range("I132").Value = 2023.125
$ws.Range("K132").Value = 6069.375
$ws.Range("M132").Value = -3539.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 421.63635
$ws.Range("I2").Value = 323.5
$ws.Range("K2").Value = 1941
$ws.Range("M2").Value = -1828
$ws.Range("H128").Value = 560030
$ws.Range("I128").Value = 560030
$ws.Range("K128").Value = 1680090
$ws.Range("M128").Value = -1675110
$ws.Range("H131").Value = 1797.5385
$ws.Range("I131").Value = 953.8
$ws.Range("J131").Value = 2324.875
$ws.Range("K131").Value = 2861.4
$ws.Range("L131").Value = 6974.625
$ws.Range("M131").Value = 2178.6
$ws.Range("N131").Value = -17054.625
$ws.Range("H137").Value = 997.5
$ws.Range("I137").Value = 997.5
$ws.Range("K137").Value = 2992.5
$ws.Range("M137").Value = 2107.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 5497.25
$ws.Range("I46").Value = 6333
$ws.Range("J46").Value = 2990
$ws.Range("K46").Value = 6333
$ws.Range("L46").Value = 2990
$ws.Range("M46").Value = -6177
$ws.Range("N46").Value = -3302
$ws.Range("H107").Value = 33334246
$ws.Range("I107").Value = 123.8
$ws.Range("K107").Value = 123.8
$ws.Range("M107").Value = 1796.2
$ws.Range("H113").Value = 617.75
$ws.Range("I113").Value = 617.75
$ws.Range("K113").Value = 617.75
$ws.Range("M113").Value = 1552.25
$ws.Range("H122").Value = 5719.5884
$ws.Range("I122").Value = 4168.25
$ws.Range("K122").Value = 12504.75
$ws.Range("M122").Value = -10054.75
$ws.Range("H126").Value = 4874.5
$ws.Range("I126").Value = 4874.5
$ws.Range("K126").Value = 14623.5
$ws.Range("M126").Value = -12153.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 8298.799999999999
$ws.Range("I35").Value = 1003
$ws.Range("J35").Value = 19242.5
$ws.Range("K35").Value = 1003
$ws.Range("L35").Value = 19242.5
$ws.Range("M35").Value = -667
$ws.Range("N35").Value = -19914.5
$ws.Range("H40").Value = 20000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 20000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 20000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -20272
$ws.Range("H46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("L46").Value = 2000
$ws.Range("N46").Value = -2376
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45997
$ws.Range("J16").Value = 45997
$ws.Range("L16").Value = 45997
$ws.Range("N16").Value = -46581
$ws.Range("H38").Value = 24000
$ws.Range("I38").Value = 24000
$ws.Range("K38").Value = 24000
$ws.Range("M38").Value = -23527
$ws.Range("H107").Value = 616.6667
$ws.Range("I107").Value = 400
$ws.Range("J107").Value = 833.3333
$ws.Range("K107").Value = 1200
$ws.Range("L107").Value = 2499.9999
$ws.Range("M107").Value = 720
$ws.Range("N107").Value = -6339.9999
$ws.Range("H113").Value = 333.75
$ws.Range("I113").Value = 280
$ws.Range("J113").Value = 495
$ws.Range("K113").Value = 840
$ws.Range("L113").Value = 1485
$ws.Range("M113").Value = 1330
$ws.Range("N113").Value = -5825
$ws.Range("H122").Value = 1999
$ws.Range("J122").Value = 1999
$ws.Range("L122").Value = 5997
$ws.Range("N122").Value = -10897
$ws.Range("H126").Value = 4239.4375
$ws.Range("I126").Value = 3394.2
$ws.Range("K126").Value = 10182.6
$ws.Range("M126").Value = -7712.599999999999
$ws.Range("H135").Value = 62542
$ws.Range("J135").Value = 62542
$ws.Range("L135").Value = 62542
$ws.Range("N135").Value = -72682
$ws.Range("H136").Value = 1295.4615
$ws.Range("I136").Value = 1295.4615
$ws.Range("K136").Value = 3886.3845
$ws.Range("M136").Value = -1336.3845
